$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three backup codes near the top of the list
$ws.Range("A2").Value = "TK50GFCXFHCN"
$ws.Range("A3").Value = "9M1A883VTX21"
$ws.Range("A4").Value = "D2ANGGG71FGC"

# Remove the two codes that were "consumed" from later in the list
$ws.Range("A15").ClearContents()
$ws.Range("A16").ClearContents()

# Update the active selection to match the saved view
$ws.Range("A6").Select()
